$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "A233"
$ws.Range("A3").Value = "A244"
$ws.Range("A4").Value = "A253"
$ws.Range("A5").Value = "A262"

$ws.Columns("A").ColumnWidth = 8.25

$ws.Range("D9").Select() | Out-Null
